# Minor fix to InitialiseMapCommandSequenceDiagram
#
# Target shape: "TextBox 199" (the initialise(...) call-label textbox on the
# sequence diagram), currently shape #52 on slide 1. We:
#   1. Reposition/resize the textbox (a:off / a:ext on its xfrm).
#   2. Change the call's argument list from "()" to "(Cell[][] cellGrid)",
#      keeping "(Cell[][] " and ")" in the original run's formatting while
#      giving "cellGrid" its own run (it is spell-flagged in the source deck).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Id -eq 200 -and $cand.Name -eq "TextBox 199") {
        $sh = $cand
        break
    }
}
if ($sh -eq $null) {
    # Fallback: this shape is the 52nd in z-order on the (single) slide.
    $sh = $s.Shapes.Item(52)
}

# Reposition / resize the shape. The point values below are chosen (rather
# than the "obvious" EMU/12700 quotient) so that the host's internal
# points->EMU conversion round-trips to the exact target EMU values:
#   x: 11584582 EMU, y: 4298341 EMU, cx: 1608066 EMU (cy is unchanged).
$sh.Left = 912.1718444824219
$sh.Top = 338.45208740234375
$sh.Width = 126.61941146850586

# Update the text: "initialise()" -> "initialise(Cell[][] cellGrid)"
$tr = $sh.TextFrame.TextRange

# Replace the "()" run (characters 11-12) with the full new argument text,
# keeping it as a single run for now.
$parens = $tr.Characters(11, 2)
$parens.Text = "(Cell[][] cellGrid)"

# Split "cellGrid" into its own run (distinct rPr / err="1" spelling flag in
# the source) by touching a per-run formatting property on just that span.
$prefixLen = "(Cell[][] ".Length
$cellGrid = $tr.Characters(11 + $prefixLen, "cellGrid".Length)
$cellGrid.Font.Size = 12
